$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 516272
$ws.Range("J17").Value = 516272
$ws.Range("L17").Value = 1548816
$ws.Range("N17").Value = -1549152
$ws.Range("H40").Value = 1765.6666
$ws.Range("I40").Value = 1849.5
$ws.Range("J40").Value = 1723.75
$ws.Range("K40").Value = 1849.5
$ws.Range("L40").Value = 1723.75
$ws.Range("M40").Value = -1674.5
$ws.Range("N40").Value = -2073.75
$ws.Range("H43").Value = 7395.7334
$ws.Range("J43").Value = 786.3333
$ws.Range("L43").Value = 786.3333
$ws.Range("N43").Value = -924.3333
$ws.Range("H138").Value = 2903.5625
$ws.Range("J138").Value = 3932.875
$ws.Range("L138").Value = 11798.625
$ws.Range("N138").Value = -22078.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4997
$ws.Range("I2").Value = 4572.4614
$ws.Range("K2").Value = 4572.4614
$ws.Range("M2").Value = -4459.4614
$ws.Range("H61").Value = 4854.04
$ws.Range("I61").Value = 4579.95
$ws.Range("J61").Value = 5950.4
$ws.Range("K61").Value = 4579.95
$ws.Range("L61").Value = 5950.4
$ws.Range("M61").Value = -4367.95
$ws.Range("N61").Value = -6374.4
$ws.Range("H74").Value = 8369.25
$ws.Range("I74").Value = 1764.8182
$ws.Range("K74").Value = 1764.8182
$ws.Range("M74").Value = -890.8181999999999
$ws.Range("H77").Value = 8369.25
$ws.Range("I77").Value = 1764.8182
$ws.Range("K77").Value = 8824.091
$ws.Range("M77").Value = -4456.091
$ws.Range("H116").Value = 4997
$ws.Range("I116").Value = 4572.4614
$ws.Range("K116").Value = 4572.4614
$ws.Range("M116").Value = -2278.4614
$ws.Range("H122").Value = 1476.0588
$ws.Range("I122").Value = 1476.0588
$ws.Range("K122").Value = 4428.1764
$ws.Range("M122").Value = -1978.1764
$ws.Range("H136").Value = 4854.04
$ws.Range("I136").Value = 4579.95
$ws.Range("J136").Value = 5950.4
$ws.Range("K136").Value = 13739.85
$ws.Range("L136").Value = 17851.2
$ws.Range("M136").Value = -11189.85
$ws.Range("N136").Value = -22951.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4997
$ws.Range("I3").Value = 4572.4614
$ws.Range("K3").Value = 4572.4614
$ws.Range("M3").Value = -4458.4614
$ws.Range("H62").Value = 65000
$ws.Range("J62").Value = 65000
$ws.Range("L62").Value = 65000
$ws.Range("N62").Value = -66372
$ws.Range("H65").Value = 65000
$ws.Range("J65").Value = 65000
$ws.Range("L65").Value = 195000
$ws.Range("N65").Value = -201864
$ws.Range("H134").Value = 2699.484
$ws.Range("I134").Value = 2353.2307
$ws.Range("K134").Value = 7059.6921
$ws.Range("M134").Value = -4524.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4547976
$ws.Range("I68").Value = 1779.8
$ws.Range("J68").Value = 8336472.5
$ws.Range("K68").Value = 5339.4
$ws.Range("L68").Value = 25009417.5
$ws.Range("M68").Value = -4528.4
$ws.Range("N68").Value = -25011039.5
$ws.Range("H71").Value = 4547976
$ws.Range("I71").Value = 1779.8
$ws.Range("J71").Value = 8336472.5
$ws.Range("K71").Value = 16018.2
$ws.Range("L71").Value = 75028252.5
$ws.Range("M71").Value = -11962.2
$ws.Range("N71").Value = -75036364.5
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("H94").Value = 4709.8887
$ws.Range("I94").Value = 2709.6
$ws.Range("K94").Value = 8128.799999999999
$ws.Range("M94").Value = -7452.799999999999
$ws.Range("H138").Value = 26325008
$ws.Range("I138").Value = 62507092
$ws.Range("J138").Value = 10765.728
$ws.Range("K138").Value = 187521276
$ws.Range("L138").Value = 32297.184
$ws.Range("M138").Value = -187516136
$ws.Range("N138").Value = -42577.18399999999
$ws.Range("N86").ClearContents()
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6400.1665
$ws.Range("I80").Value = 3726
$ws.Range("J80").Value = 11748.5
$ws.Range("K80").Value = 3726
$ws.Range("L80").Value = 11748.5
$ws.Range("M80").Value = -2728
$ws.Range("N80").Value = -13744.5
$ws.Range("H83").Value = 6400.1665
$ws.Range("I83").Value = 3726
$ws.Range("J83").Value = 11748.5
$ws.Range("K83").Value = 18630
$ws.Range("L83").Value = 58742.5
$ws.Range("M83").Value = -13638
$ws.Range("N83").Value = -68726.5
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("H113").Value = 2212.8572
$ws.Range("I113").Value = 2081.6667
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2081.6667
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 88.33329999999978
$ws.Range("N113").Value = -7340
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1748.7142
$ws.Range("I46").Value = 1454.2
$ws.Range("J46").Value = 2485
$ws.Range("K46").Value = 1454.2
$ws.Range("L46").Value = 2485
$ws.Range("M46").Value = -1266.2
$ws.Range("N46").Value = -2861
$ws.Range("H61").Value = 4691.846
$ws.Range("J61").Value = 2997.5
$ws.Range("L61").Value = 2997.5
$ws.Range("N61").Value = -3401.5
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H113").Value = 4691.846
$ws.Range("J113").Value = 2997.5
$ws.Range("L113").Value = 2997.5
$ws.Range("N113").Value = -7337.5
$ws.Range("H122").Value = 4862.6924
$ws.Range("I122").Value = 4191.5
$ws.Range("J122").Value = 7100
$ws.Range("K122").Value = 12574.5
$ws.Range("L122").Value = 21300
$ws.Range("M122").Value = -10124.5
$ws.Range("N122").Value = -26200
$ws.Range("H136").Value = 4866.4
$ws.Range("I136").Value = 4333.05
$ws.Range("J136").Value = 6999.8
$ws.Range("K136").Value = 12999.15
$ws.Range("L136").Value = 20999.4
$ws.Range("M136").Value = -10449.15
$ws.Range("N136").Value = -26099.4
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 672.0526
$ws.Range("I113").Value = 607.5454999999999
$ws.Range("J113").Value = 760.75
$ws.Range("K113").Value = 1822.6365
$ws.Range("L113").Value = 2282.25
$ws.Range("M113").Value = 347.3635000000002
$ws.Range("N113").Value = -6622.25
$ws.Range("H132").Value = 3125.4082
$ws.Range("I132").Value = 3040
$ws.Range("K132").Value = 9120
$ws.Range("M132").Value = -6590
$ws.Range("H136").Value = 1782.625
$ws.Range("I136").Value = 1744.2222
$ws.Range("J136").Value = 1990
$ws.Range("K136").Value = 5232.6666
$ws.Range("L136").Value = 5970
$ws.Range("M136").Value = -2682.6666
